$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from E1 to F1, then set header text
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("F1").Value = "time_taken"

# Fill in time_taken values for each data row
$ws.Range("F2").Value = "2021-10-05 13:41:46.264382"
$ws.Range("F3").Value = "2021-10-05 13:41:46.264393"
$ws.Range("F4").Value = "2021-10-05 13:41:46.264396"
$ws.Range("F5").Value = "2021-10-05 13:41:46.264399"
$ws.Range("F6").Value = "2021-10-05 13:41:46.264402"
$ws.Range("F7").Value = "2021-10-05 13:41:46.264405"
$ws.Range("F8").Value = "2021-10-05 13:41:46.264408"
$ws.Range("F9").Value = "2021-10-05 13:41:46.264410"
$ws.Range("F10").Value = "2021-10-05 13:41:46.264413"
$ws.Range("F11").Value = "2021-10-05 13:41:46.264416"
$ws.Range("F12").Value = "2021-10-05 13:41:46.264419"
$ws.Range("F13").Value = "2021-10-05 13:41:46.264421"
$ws.Range("F14").Value = "2021-10-05 13:41:46.264424"
$ws.Range("F15").Value = "2021-10-05 13:41:46.264426"
$ws.Range("F16").Value = "2021-10-05 13:41:46.264429"
$ws.Range("F17").Value = "2021-10-05 13:41:46.264432"
$ws.Range("F18").Value = "2021-10-05 13:41:46.264435"
$ws.Range("F19").Value = "2021-10-05 13:41:46.264438"
$ws.Range("F20").Value = "2021-10-05 13:41:46.264441"
$ws.Range("F21").Value = "2021-10-05 13:41:46.264443"
$ws.Range("F22").Value = "2021-10-05 13:41:46.264446"
$ws.Range("F23").Value = "2021-10-05 13:41:46.264448"
$ws.Range("F24").Value = "2021-10-05 13:41:46.264451"
$ws.Range("F25").Value = "2021-10-05 13:41:46.264454"
$ws.Range("F26").Value = "2021-10-05 13:41:46.264457"
$ws.Range("F27").Value = "2021-10-05 13:41:46.264460"
$ws.Range("F28").Value = "2021-10-05 13:41:46.264462"
$ws.Range("F29").Value = "2021-10-05 13:41:46.264465"
$ws.Range("F30").Value = "2021-10-05 13:41:46.264468"
$ws.Range("F31").Value = "2021-10-05 13:41:46.264470"
$ws.Range("F32").Value = "2021-10-05 13:41:46.264473"
$ws.Range("F33").Value = "2021-10-05 13:41:46.264476"
$ws.Range("F34").Value = "2021-10-05 13:41:46.264479"
$ws.Range("F35").Value = "2021-10-05 13:41:46.264482"
$ws.Range("F36").Value = "2021-10-05 13:41:46.264485"
$ws.Range("F37").Value = "2021-10-05 13:41:46.264487"
$ws.Range("F38").Value = "2021-10-05 13:41:46.264490"
$ws.Range("F39").Value = "2021-10-05 13:41:46.264493"
$ws.Range("F40").Value = "2021-10-05 13:41:46.264495"
$ws.Range("F41").Value = "2021-10-05 13:41:46.264498"
$ws.Range("F42").Value = "2021-10-05 13:41:46.264501"
$ws.Range("F43").Value = "2021-10-05 13:41:46.264504"
$ws.Range("F44").Value = "2021-10-05 13:41:46.264506"
$ws.Range("F45").Value = "2021-10-05 13:41:46.264509"
$ws.Range("F46").Value = "2021-10-05 13:41:46.264512"
$ws.Range("F47").Value = "2021-10-05 13:41:46.264515"
$ws.Range("F48").Value = "2021-10-05 13:41:46.264518"
$ws.Range("F49").Value = "2021-10-05 13:41:46.264520"
$ws.Range("F50").Value = "2021-10-05 13:41:46.264523"
$ws.Range("F51").Value = "2021-10-05 13:41:46.264526"
$ws.Range("F52").Value = "2021-10-05 13:41:46.264529"
$ws.Range("F53").Value = "2021-10-05 13:41:46.264531"
$ws.Range("F54").Value = "2021-10-05 13:41:46.264535"
$ws.Range("F55").Value = "2021-10-05 13:41:46.264538"
$ws.Range("F56").Value = "2021-10-05 13:41:46.264541"
$ws.Range("F57").Value = "2021-10-05 13:41:46.264544"
